# Apply the Dec 22 2022 symbol-list refresh (GitHub Actions) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Sheet, $Ref, $NewValue) {
    # Write $NewValue into $Ref as TEXT (not auto-coerced to a number),
    # then restore the cell style so no formatting diff is introduced.
    $cell = $Sheet.Range($Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "245.79"
Set-TextCell $ws "D3" "22.56"
Set-TextCell $ws "D4" "5.392"
Set-TextCell $ws "D5" "0.05760"
Set-TextCell $ws "D7" "6.323"
Set-TextCell $ws "D8" "0.8105"
Set-TextCell $ws "D9" "0.8929"
Set-TextCell $ws "D10" "0.1441"
Set-TextCell $ws "D11" "0.07330"
Set-TextCell $ws "D12" "0.03128"
Set-TextCell $ws "E12" "11BitrueCoinBTRBestin24h"
Set-TextCell $ws "D13" "0.02973"
Set-TextCell $ws "D14" "0.09415"
Set-TextCell $ws "B15" "MCDex"
Set-TextCell $ws "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell $ws "D15" "3.932"
Set-TextCell $ws "E15" "14MCDexMCB"
Set-TextCell $ws "B16" "BitForexToken"
Set-TextCell $ws "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell $ws "D16" "0.001583"
Set-TextCell $ws "E16" "15BitForexTokenBF"
Set-TextCell $ws "B17" "CoinExToken"
Set-TextCell $ws "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell $ws "D17" "0.04795"
Set-TextCell $ws "E17" "16CoinExTokenCET"
Set-TextCell $ws "B18" "One"
Set-TextCell $ws "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell $ws "D18" "0.0005853"
Set-TextCell $ws "E18" "17OneONE"
Set-TextCell $ws "B19" "TigerCash"
Set-TextCell $ws "C19" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell $ws "D19" "0.006355"
Set-TextCell $ws "E19" "18TigerCashTCH"
Set-TextCell $ws "B20" "HotbitToken"
Set-TextCell $ws "C20" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell $ws "D20" "0.004063"
Set-TextCell $ws "E20" "19HotbitTokenHTB"
Set-TextCell $ws "B21" "BitKan"
Set-TextCell $ws "C21" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell $ws "D21" "0.0009927"
Set-TextCell $ws "E21" "20BitKanKAN"
Set-TextCell $ws "B22" "NitroEx"
Set-TextCell $ws "C22" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextCell $ws "D22" "0.0001500"
Set-TextCell $ws "E22" "21NitroExNTX"
Set-TextCell $ws "B23" "LEO"
Set-TextCell $ws "C23" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D23" "3.721"
Set-TextCell $ws "E23" "22LEOLEO"
Set-TextCell $ws "B24" "BTSEToken"
Set-TextCell $ws "C24" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell $ws "D24" "2.193"
Set-TextCell $ws "E24" "23BTSETokenBTSE"
Set-TextCell $ws "B25" "BitpandaEcosystemToken"
Set-TextCell $ws "C25" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell $ws "D25" "0.3272"
Set-TextCell $ws "E25" "24BitpandaEcosystemTokenBEST"
Set-TextCell $ws "B26" "ProBitToken"
Set-TextCell $ws "C26" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell $ws "D26" "0.1301"
Set-TextCell $ws "E26" "25ProBitTokenPROB"
Set-TextCell $ws "D40" "0.03896"
Set-TextCell $ws "D41" "0.006791"
Set-TextCell $ws "D43" "0.002411"
Set-TextCell $ws "E43" "42CEJICEJI"
Set-TextCell $ws "D45" "0.00005651"
Set-TextCell $ws "D48" "0.1647"
Set-TextCell $ws "D50" "0.01010"

Write-Output "Applied 68 cell updates"
